$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12, pushing existing rows 12-29 down to 13-30.
$ws.Rows.Item(12).Insert()

# Copy the style (date format) from D13 (the row that was pushed down, formerly D12) into the new D12.
$ws.Range("D13").Copy()
$ws.Range("D12").PasteSpecial(-4122) | Out-Null

# Populate the new row 12 with data.
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").Value = 45272
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100103
$ws.Range("H12").Value = "Frutos de hueso (carozo)"
$ws.Range("I12").Value = 100103001
$ws.Range("J12").Value = "Cereza"
$ws.Range("K12").Value = "Lapins"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 230
$ws.Range("N12").Value = 23000
$ws.Range("O12").Value = 25000
$ws.Range("P12").Value = 24304
$ws.Range("Q12").Value = "`$/caja 20 kilos"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 1215
$ws.Range("T12").Value = 20
